$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New test data (rows 2-5), keyed by current column layout:
# A=firstname B=lastname C=email(to be removed) D=phone E=password F=occupation G=gender
$firstnames = @("testuser5", "testuser6", "testuser7", "testuser8")
$lastnames  = @("lastname5", "lastname6", "lastname7", "lastname8")
$phones     = @("1234567898", "2234567898", "3123456789", "4234567898")

for ($i = 0; $i -lt 4; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $firstnames[$i]
    $ws.Cells.Item($row, 2).Value = $lastnames[$i]
    # Store the phone number as text (quote-prefixed), same as column D before the
    # email column is removed; after removal this becomes column C.
    $ws.Cells.Item($row, 4).Value = "'" + $phones[$i]
}

# Remove the now-unused "email" column entirely; everything to its right
# (phone, password, occupation, gender) shifts one column to the left.
$ws.Range("C1").EntireColumn.Delete()
